# Add the "Metadata" worksheet (second sheet, after Sheet1) and populate it
# with a data-dictionary describing the columns of Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Change Sheet1's selection before we touch the new sheet, so the new sheet
# ends up as the active / selected tab afterwards.
$ws1.Range("A1:H1").Select() | Out-Null

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Metadata"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 9
$ws2.Columns.Item(2).ColumnWidth = 77.75

# Data dictionary rows: column name in A, description in B
$ws2.Range("A2").Value = "Block"
$ws2.Range("B2").Value = "5 blocks in the Randomized Complete Block design to capture field variablity"

$ws2.Range("A3").Value = "Trt"
$ws2.Range("B3").WrapText = $true
$ws2.Range("B3").Value = "Treatments: 1 - Abound low rate, 2 - Abound high rate, 3 - Sovran low rate, 4- Sovran high rate, 5 standard fungicide(pos. control), 6 - Water check (neg. control)"
$ws2.Rows.Item(3).RowHeight = 28

$ws2.Range("A4").Value = "LeafInc"
$ws2.Range("B4").Value = "Percent of leaves with phomopsis infection"

$ws2.Range("A6").Value = "NodeInc"
$ws2.Range("B6").Value = "Percent of nodes with phomopsis infection"

$ws2.Range("A8").Value = "ClusterInc"
$ws2.Range("B8").Value = "Percent of clusters with phomopsis infection"

$ws2.Range("A5").Value = "LeafSev"
$ws2.Range("B5").Value = "Percent of area of leaves with phomopsis infection determined from Barratt Horsfall Scale"

$ws2.Range("A7").Value = "NodeSev"
$ws2.Range("B7").Value = "Percent of area of nodes  with phomopsis infection  determined from Barratt Horsfall Scale"

$ws2.Range("A9").Value = "ClusterSev"
$ws2.Range("B9").Value = "Percent of area of clusters with phomopsis infection  determined from Barratt Horsfall Scale"

# Title cell (added last so the "Metadata" shared string sorts after the
# other new strings, matching the authored sharedStrings order)
$ws2.Range("B1").Value = "Metadata"
$ws2.Range("B1").Font.Bold = $true

$ws2.Range("B1").Select() | Out-Null
